$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value2 = 106002741
$ws.Range("B3").Value2 = 98520
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value2 = "Ovaliderad"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "LC"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = 222498
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "Blåsippa"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value2 = "Hepatica nobilis"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value2 = "Schreb."
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value2 = $null
$ws.Range("J3").Value2 = $null
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value2 = "Bergstorp , Srm"
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Value2 = 580759.2328223517
$ws.Range("R3").Value2 = 6571012.762967644
$ws.Range("S3").Value2 = 4
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value2 = "Södermanland"
$ws.Range("T3").Style = "Normal"
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value2 = "Eskilstuna"
$ws.Range("U3").Style = "Normal"
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value2 = "Södermanland"
$ws.Range("V3").Style = "Normal"
$ws.Range("W3").NumberFormat = "@"
$ws.Range("W3").Value2 = "Husby-Rekarne"
$ws.Range("W3").Style = "Normal"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value2 = "2023-01-15"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value2 = "12:41"
$ws.Range("Z3").Style = "Normal"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value2 = "2023-01-15"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value2 = "12:41"
$ws.Range("AB3").Style = "Normal"
$ws.Range("AC3").Value2 = $null
$ws.Range("AD3").Value2 = $false
$ws.Range("AE3").Value2 = $false
$ws.Range("AG3").Value2 = $false
$ws.Range("AW3").NumberFormat = "@"
$ws.Range("AW3").Value2 = "Michael Lander"
$ws.Range("AW3").Style = "Normal"
$ws.Range("AX3").NumberFormat = "@"
$ws.Range("AX3").Value2 = "Michael Lander"
$ws.Range("AX3").Style = "Normal"

# Row 4
$ws.Range("A4").Value2 = 106538494
$ws.Range("B4").Value2 = 103265
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value2 = "Ovaliderad"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "LC"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = 221144
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value2 = "Grönpyrola"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value2 = "Pyrola chlorantha"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value2 = "Sw."
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value2 = "10"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value2 = "plantor/tuvor"
$ws.Range("J4").Style = "Normal"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value2 = "Bergstorp , Srm"
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value2 = 580735.7162288441
$ws.Range("R4").Value2 = 6571227.517849359
$ws.Range("S4").Value2 = 4
$ws.Range("T4").NumberFormat = "@"
$ws.Range("T4").Value2 = "Södermanland"
$ws.Range("T4").Style = "Normal"
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value2 = "Eskilstuna"
$ws.Range("U4").Style = "Normal"
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value2 = "Södermanland"
$ws.Range("V4").Style = "Normal"
$ws.Range("W4").NumberFormat = "@"
$ws.Range("W4").Value2 = "Husby-Rekarne"
$ws.Range("W4").Style = "Normal"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2023-02-08"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value2 = "11:19"
$ws.Range("Z4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2023-02-08"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value2 = "11:19"
$ws.Range("AB4").Style = "Normal"
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value2 = "Ett tiotal plantor inom ett område av 4 kvm"
$ws.Range("AC4").Style = "Normal"
$ws.Range("AD4").Value2 = $false
$ws.Range("AE4").Value2 = $false
$ws.Range("AG4").Value2 = $false
$ws.Range("AW4").NumberFormat = "@"
$ws.Range("AW4").Value2 = "Michael Lander"
$ws.Range("AW4").Style = "Normal"
$ws.Range("AX4").NumberFormat = "@"
$ws.Range("AX4").Value2 = "Michael Lander"
$ws.Range("AX4").Style = "Normal"
